$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two oldest "Periodo Mora" entries (2504) for the existing
# workers - previous EC rows are no longer needed. Deleting rows 16:17
# shifts the remaining four data rows (old 18-21, period 2505/2506) up
# into rows 16-19, carrying the bottom-border style of the old last row
# (21) onto the new last row (19) automatically.
$ws.Rows("16:17").Delete()

# New worker 1: FREDY DE JESUS SIERRA VARELA
$ws.Range("C16").Value = "73167275"
$ws.Range("D16").Value = "FREDY DE JESUS SIERRA VARELA"
$ws.Range("E16").Value = "2507"
$ws.Range("F16").Value = 120000
$ws.Range("G16").Value = 3000000

# New worker 2: IVAN ARTURO BROCHET BAYONA
$ws.Range("C17").Value = "73139008"
$ws.Range("D17").Value = "IVAN ARTURO BROCHET BAYONA"
$ws.Range("E17").Value = "2507"
$ws.Range("F17").Value = 56940
$ws.Range("G17").Value = 1423500

# Existing worker DANIELA FLOREZ CAMARGO, moved to period 2507
$ws.Range("C18").Value = "1047463262"
$ws.Range("D18").Value = "DANIELA FLOREZ CAMARGO"
$ws.Range("E18").Value = "2507"
$ws.Range("F18").Value = 57520
$ws.Range("G18").Value = 1438000

# Existing worker KATYA CAROLINA HAWKINS RAMIREZ, moved to period 2507
$ws.Range("C19").Value = "1143164560"
$ws.Range("D19").Value = "KATYA CAROLINA HAWKINS RAMIREZ"
$ws.Range("E19").Value = "2507"
$ws.Range("F19").Value = 56940
$ws.Range("G19").Value = 1423500

# Update the summary figures: total overdue value, worker count, period count
$ws.Range("E11").Value = 291400
$ws.Range("C13").Value = 4
$ws.Range("F13").Value = 1
